# Updated symbol list on Sun Dec 11 23:57:26 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) refreshes -------------------------------------------------
# These cells hold numeric-looking text values (t="inlineStr" in the source file),
# so force the cell format to Text before assigning, otherwise Excel would store
# them as plain numbers and drop things like trailing zeros.
$priceUpdates = @{
    "D3"  = "21.09"
    "D4"  = "6.450"
    "D7"  = "1.531"
    "D8"  = "6.563"
    "D9"  = "0.8200"
    "D12" = "0.08630"
    "D14" = "0.03214"
    "D15" = "0.09196"
    "D16" = "3.708"
    "D17" = "0.001650"
    "D18" = "0.04764"
    "D19" = "0.006133"
    "D20" = "0.006271"
    "D23" = "3.784"
    "D25" = "0.3355"
    "D26" = "0.1261"
    "D40" = "0.04763"
    "D41" = "0.007129"
    "D42" = "0.1116"
    "D43" = "0.003454"
    "D44" = "0.01146"
    "D45" = "0.00006903"
    "D47" = "0.9010"
    "D48" = "0.003117"
    "D50" = "0.01241"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# --- Row 42 and row 43 swap their coin/link identity (symbol list reordered) ----
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
